# Applies scheduled-runner market-price/profit updates to the Leve profit
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), cell by cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1704.3334
$ws.Range("I20").Value = 1704.3334
$ws.Range("K20").Value = 1704.3334
$ws.Range("M20").Value = -1474.3334
$ws.Range("H32").Value = 957.7143
$ws.Range("I32").Value = 767
$ws.Range("J32").Value = 1100.75
$ws.Range("K32").Value = 767
$ws.Range("L32").Value = 1100.75
$ws.Range("M32").Value = -441
$ws.Range("N32").Value = -1752.75
$ws.Range("H35").Value = 1704.3334
$ws.Range("I35").Value = 1704.3334
$ws.Range("K35").Value = 1704.3334
$ws.Range("M35").Value = -1325.3334
$ws.Range("H116").Value = 1560.5454
$ws.Range("I116").Value = 1272.8572
$ws.Range("J116").Value = 2064
$ws.Range("K116").Value = 1272.8572
$ws.Range("L116").Value = 2064
$ws.Range("M116").Value = 2169.1428
$ws.Range("N116").Value = -8948
$ws.Range("H138").Value = 2647991.8
$ws.Range("I138").Value = 1679.4615
$ws.Range("J138").Value = 3336033
$ws.Range("K138").Value = 5038.3845
$ws.Range("L138").Value = 10008099
$ws.Range("M138").Value = 101.6154999999999
$ws.Range("N138").Value = -10018379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1604.881
$ws.Range("I2").Value = 1304.8387
$ws.Range("J2").Value = 2450.4546
$ws.Range("K2").Value = 1304.8387
$ws.Range("L2").Value = 2450.4546
$ws.Range("M2").Value = -1191.8387
$ws.Range("N2").Value = -2676.4546
$ws.Range("H61").Value = 40082084
$ws.Range("I61").Value = 83418260
$ws.Range("K61").Value = 83418260
$ws.Range("M61").Value = -83418048
$ws.Range("H63").Value = 3042.4614
$ws.Range("I63").Value = 2666.2856
$ws.Range("K63").Value = 2666.2856
$ws.Range("M63").Value = -1980.2856
$ws.Range("H66").Value = 3042.4614
$ws.Range("I66").Value = 2666.2856
$ws.Range("K66").Value = 13331.428
$ws.Range("M66").Value = -9899.428
$ws.Range("H74").Value = 7413788
$ws.Range("I74").Value = 11954366
$ws.Range("J74").Value = 79007.69500000001
$ws.Range("K74").Value = 11954366
$ws.Range("L74").Value = 79007.69500000001
$ws.Range("M74").Value = -11953492
$ws.Range("N74").Value = -80755.69500000001
$ws.Range("H77").Value = 7413788
$ws.Range("I77").Value = 11954366
$ws.Range("J77").Value = 79007.69500000001
$ws.Range("K77").Value = 59771830
$ws.Range("L77").Value = 395038.475
$ws.Range("M77").Value = -59767462
$ws.Range("N77").Value = -403774.475
$ws.Range("H110").Value = 2017.8889
$ws.Range("I110").Value = 1005.5
$ws.Range("J110").Value = 2827.8
$ws.Range("K110").Value = 1005.5
$ws.Range("L110").Value = 2827.8
$ws.Range("M110").Value = 1039.5
$ws.Range("N110").Value = -6917.8
$ws.Range("H116").Value = 1604.881
$ws.Range("I116").Value = 1304.8387
$ws.Range("J116").Value = 2450.4546
$ws.Range("K116").Value = 1304.8387
$ws.Range("L116").Value = 2450.4546
$ws.Range("M116").Value = 989.1613
$ws.Range("N116").Value = -7038.4546
$ws.Range("H122").Value = 7938739.5
$ws.Range("I122").Value = 2204
$ws.Range("J122").Value = 37039372
$ws.Range("K122").Value = 6612
$ws.Range("L122").Value = 111118116
$ws.Range("M122").Value = -4162
$ws.Range("N122").Value = -111123016
$ws.Range("H132").Value = 72034.8
$ws.Range("I132").Value = 72978.86
$ws.Range("J132").Value = 71208.75
$ws.Range("K132").Value = 218936.58
$ws.Range("L132").Value = 213626.25
$ws.Range("M132").Value = -216406.58
$ws.Range("N132").Value = -218686.25
$ws.Range("H136").Value = 40082084
$ws.Range("I136").Value = 83418260
$ws.Range("K136").Value = 250254780
$ws.Range("M136").Value = -250252230
$ws.Range("H138").Value = 40318.57
$ws.Range("J138").Value = 40318.57
$ws.Range("L138").Value = 40318.57
$ws.Range("N138").Value = -50598.57
$ws.Range("H141").Value = 46000
$ws.Range("J141").Value = 46000
$ws.Range("L141").Value = 46000
$ws.Range("N141").Value = -56360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1604.881
$ws.Range("I3").Value = 1304.8387
$ws.Range("J3").Value = 2450.4546
$ws.Range("K3").Value = 1304.8387
$ws.Range("L3").Value = 2450.4546
$ws.Range("M3").Value = -1190.8387
$ws.Range("N3").Value = -2678.4546
$ws.Range("H20").Value = 1022.619
$ws.Range("I20").Value = 816.5
$ws.Range("J20").Value = 1297.4445
$ws.Range("K20").Value = 816.5
$ws.Range("L20").Value = 1297.4445
$ws.Range("M20").Value = -569.5
$ws.Range("N20").Value = -1791.4445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 62900.125
$ws.Range("I10").Value = 71600.14
$ws.Range("K10").Value = 71600.14
$ws.Range("M10").Value = -71461.14
$ws.Range("H31").Value = 2888.3215
$ws.Range("I31").Value = 1223.65
$ws.Range("J31").Value = 7050
$ws.Range("K31").Value = 1223.65
$ws.Range("L31").Value = 7050
$ws.Range("M31").Value = -928.6500000000001
$ws.Range("N31").Value = -7640
$ws.Range("H34").Value = 2888.3215
$ws.Range("I34").Value = 1223.65
$ws.Range("J34").Value = 7050
$ws.Range("K34").Value = 1223.65
$ws.Range("L34").Value = 7050
$ws.Range("M34").Value = -1021.65
$ws.Range("N34").Value = -7454
$ws.Range("H58").Value = 47622492
$ws.Range("I58").Value = 66669564
$ws.Range("J58").Value = 4816.8335
$ws.Range("K58").Value = 66669564
$ws.Range("L58").Value = 4816.8335
$ws.Range("M58").Value = -66669361
$ws.Range("N58").Value = -5222.8335
$ws.Range("H132").Value = 66492.06
$ws.Range("I132").Value = 4441.727
$ws.Range("J132").Value = 203002.8
$ws.Range("K132").Value = 13325.181
$ws.Range("L132").Value = 609008.3999999999
$ws.Range("M132").Value = -10795.181
$ws.Range("N132").Value = -614068.3999999999
$ws.Range("H134").Value = 47411.043
$ws.Range("I134").Value = 3294.5715
$ws.Range("J134").Value = 109174.1
$ws.Range("K134").Value = 9883.7145
$ws.Range("L134").Value = 327522.3
$ws.Range("M134").Value = -7348.7145
$ws.Range("N134").Value = -332592.3
$ws.Range("H136").Value = 47622492
$ws.Range("I136").Value = 66669564
$ws.Range("J136").Value = 4816.8335
$ws.Range("K136").Value = 200008692
$ws.Range("L136").Value = 14450.5005
$ws.Range("M136").Value = -200006142
$ws.Range("N136").Value = -19550.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 989.8
$ws.Range("I98").Value = 999.6667
$ws.Range("J98").Value = 975
$ws.Range("K98").Value = 2999.0001
$ws.Range("L98").Value = 2925
$ws.Range("M98").Value = -1501.0001
$ws.Range("N98").Value = -5921
$ws.Range("H131").Value = 951.5733
$ws.Range("J131").Value = 1002.1061
$ws.Range("L131").Value = 3006.3183
$ws.Range("N131").Value = -13086.3183

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1071.6666
$ws.Range("I102").Value = 1071.6666
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1071.6666
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 550.3334
$ws.Range("N102").Value = $null
$ws.Range("H126").Value = 2329.2632
$ws.Range("I126").Value = 1314.2858
$ws.Range("J126").Value = 2921.3333
$ws.Range("K126").Value = 3942.8574
$ws.Range("L126").Value = 8763.999899999999
$ws.Range("M126").Value = -1472.8574
$ws.Range("N126").Value = -13703.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 295.85715
$ws.Range("I55").Value = 216.33333
$ws.Range("J55").Value = 355.5
$ws.Range("K55").Value = 216.33333
$ws.Range("L55").Value = 355.5
$ws.Range("M55").Value = -43.33332999999999
$ws.Range("N55").Value = -701.5
$ws.Range("H100").Value = 49036.953
$ws.Range("I100").Value = 72555.42999999999
$ws.Range("K100").Value = 72555.42999999999
$ws.Range("M100").Value = -72014.42999999999
$ws.Range("H136").Value = 251024.88
$ws.Range("I136").Value = 167533.33
$ws.Range("K136").Value = 502599.99
$ws.Range("M136").Value = -500049.99

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 272.31033
$ws.Range("I107").Value = 302.53333
$ws.Range("J107").Value = 239.92857
$ws.Range("K107").Value = 907.5999899999999
$ws.Range("L107").Value = 719.78571
$ws.Range("M107").Value = 1012.40001
$ws.Range("N107").Value = -4559.78571
$ws.Range("H136").Value = 47552.906
$ws.Range("I136").Value = 24875.596
$ws.Range("K136").Value = 74626.788
$ws.Range("M136").Value = -72076.788
